$d = $word.ActiveDocument
$bullet = [char]0x2022

# ------------------------------------------------------------------
# 1. Contact-info line: merge "github.com/" + "danielmartincraig" (which
#    were split across runs with a proofed spell-check range) into one
#    run, and append the previously-separate " linkedin.com/..." part.
#    Doing the Find/Replace across the whole line collapses every run
#    in the match into a single plain run and drops the stale
#    <w:proofErr/> markers.
# ------------------------------------------------------------------
$contactOld = "(803)389-6750 $bullet danielmartincraig@gmail.com $bullet github.com/danielmartincraig $bullet linkedin.com/danielcraig23"
$contactRange = $d.Content
$found = $contactRange.Find.Execute($contactOld, $true, $false, $false, $false, $false, $true, 1, $false, $contactOld, 2)
$contactPara = $contactRange.Paragraphs(1)

# ------------------------------------------------------------------
# 2. Remove the old "_GoBack" bookmark that used to sit at the very
#    end of the document (near "Fluent in Spanish") -- it moves to the
#    new OBJECTIVE paragraph below. Do this BEFORE inserting the new
#    bookmark so the "_GoBack" name is unambiguous (only one instance
#    exists in the document at this point).
# ------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks("_GoBack").Delete()
}

# ------------------------------------------------------------------
# 3. Insert a new "OBJECTIVE:" heading paragraph right after the
#    contact-info paragraph.
# ------------------------------------------------------------------
$contactPara.Range.InsertParagraphAfter()
$objectivePara = $contactPara.Next()
$objectiveXml = '<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
  '<w:pPr><w:pStyle w:val="Heading1"/><w:rPr><w:sz w:val="24"/></w:rPr></w:pPr>' +
  '<w:r><w:t xml:space="preserve">OBJECTIVE: </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve">Eager to drive back-end solutions at </w:t></w:r>' +
  '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t>Dell</w:t></w:r>' +
  '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>' +
  '<w:r><w:rPr><w:sz w:val="24"/></w:rPr><w:t xml:space="preserve"> on a full-time basis</w:t></w:r>' +
  '</w:p>'
$objectivePara.Range.InsertXML($objectiveXml)

# ------------------------------------------------------------------
# 4. "Web Engineering I and II" bullet: merge the three runs (split by
#    a grammar-check proofed range around "Engineering") into one run.
# ------------------------------------------------------------------
$webOld = $bullet + "    Web Engineering I and II"
$found2 = $d.Content.Find.Execute($webOld, $true, $false, $false, $false, $false, $true, 1, $false, $webOld, 2)

Write-Output "contact found: $found; web found: $found2"
